{"js": "// Department, study program and subject management\n//\n// 1) \"\u0160kolite\u013e:\" label gains a trailing space.\n// 2) The \"{Supervisor}\" placeholder run drops its now-redundant leading space.\n// 3) The {Subject1}/{Subject2}/{Subject3} table rows collapse into a single\n//    {Subjects} row, and the (now single-row) column narrows from 1019 to\n//    994 twentieths of a point (49.7pt).\n\nconst body = context.document.body;\n\n// --- 1) \"\u0160kolite\u013e:\" -> \"\u0160kolite\u013e: \" --------------------------------------\nconst skolitelLabel = body.search(\"\u0160kolite\u013e:\", { matchCase: true });\nskolitelLabel.load(\"items\");\nawait context.sync();\nif (skolitelLabel.items.length > 0) {\n  skolitelLabel.items[0].insertText(\"\u0160kolite\u013e: \", Word.InsertLocation.replace);\n}\n\n// --- 2) \" {Supervisor}\" -> \"{Supervisor}\" --------------------------------\nconst supervisorRun = body.search(\" {Supervisor}\", { matchCase: true });\nsupervisorRun.load(\"items\");\nawait context.sync();\nif (supervisorRun.items.length > 0) {\n  supervisorRun.items[0].insertText(\"{Supervisor}\", Word.InsertLocation.replace);\n}\n\n// --- 3) Subject table: merge Subject1/Subject2/Subject3 into Subjects ----\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nif (rows.length > 0) {\n  // Rename the placeholder that stays behind, in the first row.\n  const subjectText = rows[0].search(\"{Subject1}\", { matchCase: true });\n  subjectText.load(\"items\");\n  await context.sync();\n  if (subjectText.items.length > 0) {\n    subjectText.items[0].insertText(\"{Subjects}\", Word.InsertLocation.replace);\n  }\n\n  // Drop the now-unneeded second and third rows.\n  for (let i = rows.length - 1; i >= 1; i--) {\n    rows[i].delete();\n  }\n  await context.sync();\n\n  // Narrow the remaining column from 1019 to 994 twentieths of a point.\n  table.rows.load(\"items\");\n  await context.sync();\n  const firstRow = table.rows.items[0];\n  firstRow.cells.load(\"items\");\n  await context.sync();\n  firstRow.cells.items[0].columnWidth = 49.7;\n  await context.sync();\n}\n", "ps1": "# Department, study program and subject management\n#\n# 1) \"\u0160kolite\u013e:\" label gains a trailing space.\n# 2) The \"{Supervisor}\" placeholder run drops its now-redundant leading space.\n# 3) The {Subject1}/{Subject2}/{Subject3} table rows collapse into a single\n#    {Subjects} row, and the (now single-row) column narrows from 1019 to\n#    994 twentieths of a point (49.7pt).\n\n$d = $word.ActiveDocument\n\n# --- 1) \"\u0160kolite\u013e:\" -> \"\u0160kolite\u013e: \" --------------------------------------\n$d.Content.Find.Execute(\"\u0160kolite\u013e:\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0160kolite\u013e: \", 2)\n\n# --- 2) \" {Supervisor}\" -> \"{Supervisor}\" --------------------------------\n$d.Content.Find.Execute(\" {Supervisor}\", $false, $false, $false, $false, $false, $true, 1, $false, \"{Supervisor}\", 2)\n\n# --- 3) Subject table: merge Subject1/Subject2/Subject3 into Subjects ----\n$t = $d.Tables.Item(1)\n\n# Rename the placeholder that stays behind, in the first row.\n$t.Rows.Item(1).Range.Find.Execute(\"{Subject1}\", $false, $false, $false, $false, $false, $true, 1, $false, \"{Subjects}\", 2)\n\n# Drop the now-unneeded second and third rows (delete from the bottom up).\n$t.Rows.Item(3).Delete()\n$t.Rows.Item(2).Delete()\n\n# Narrow the remaining column from 1019 to 994 twentieths of a point.\n$t.Columns.Item(1).Width = 49.7\n"}
